# Appends 7 new data rows (583-589) to the Landscaping Data log, continuing
# the existing table on Sheet1 (same columns A:T used throughout the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 583
$endRow = 589

# Column data, one flat array per column, in row order 583..589
$aVals = @(45870, 45870, 45870, 45870, 45870, 45870, 45870)
$bVals = @("Flowering", "Nonflowering", "Nonflowering", "Nonflowering", "Nonflowering", "Nonflowering", "Tree")
$cVals = @("Large", "Medium", "Small", "Medium", "Medium", "Large", "Medium")
$dVals = @(58, 58, 58, 58, 58, 58, 58)
$eVals = @(78, 78, 78, 78, 78, 78, 78)
$gVals = @(0, 0, 0, 0, 0, 0, 0)
$hVals = @(0, 0, 0.1, 0.25, 0.3, 0.5, 1.5)
$iVals = @("No", "No", "No", "No", "No", "No", "No")
$jVals = @(2, 3, 3, 3, 3, 4, 1)
$kVals = @("Bright", "Bright", "Neutral", "Neutral", "Bright", "Neutral", "Neutral")
$lVals = @(7, 7, 7, 7, 7, 7, 7)
$mVals = @(0.61, 0.61, 0.61, 0.61, 0.61, 0.61, 0.61)
$nVals = @(57, 57, 57, 57, 57, 57, 57)
$oVals = @(30.22, 30.22, 30.22, 30.22, 30.22, 30.22, 30.22)
$pVals = @(12, 12, 12, 12, 12, 12, 12)
$qVals = @(0.26, 0.26, 0.26, 0.26, 0.26, 0.26, 0.26)
$rVals = @(9.9, 9.9, 9.9, 9.9, 9.9, 9.9, 9.9)
$sVals = @(69, 69, 69, 69, 69, 69, 69)
$tVals = @(5, 5, 5, 5, 5, 5, 5)

for ($i = 0; $i -lt 7; $i++) {
  $r = $startRow + $i

  # Column A keeps the same date-number formatting as the row above it
  $ws.Cells.Item($r - 1, 1).Copy()
  $ws.Cells.Item($r, 1).PasteSpecial(-4122)
  $ws.Cells.Item($r, 1).Value = $aVals[$i]

  $ws.Cells.Item($r, 2).Value = $bVals[$i]
  $ws.Cells.Item($r, 3).Value = $cVals[$i]
  $ws.Cells.Item($r, 4).Value = $dVals[$i]
  $ws.Cells.Item($r, 5).Value = $eVals[$i]
  $ws.Cells.Item($r, 7).Value = $gVals[$i]
  $ws.Cells.Item($r, 8).Value = $hVals[$i]
  $ws.Cells.Item($r, 9).Value = $iVals[$i]
  $ws.Cells.Item($r, 10).Value = $jVals[$i]
  $ws.Cells.Item($r, 11).Value = $kVals[$i]
  $ws.Cells.Item($r, 12).Value = $lVals[$i]
  $ws.Cells.Item($r, 13).Value = $mVals[$i]
  $ws.Cells.Item($r, 14).Value = $nVals[$i]
  $ws.Cells.Item($r, 15).Value = $oVals[$i]
  $ws.Cells.Item($r, 16).Value = $pVals[$i]
  $ws.Cells.Item($r, 17).Value = $qVals[$i]
  $ws.Cells.Item($r, 18).Value = $rVals[$i]
  $ws.Cells.Item($r, 19).Value = $sVals[$i]
  $ws.Cells.Item($r, 20).Value = $tVals[$i]
}

# Fill the F (Temp_Diff) column for the new rows in one shot so it is
# emitted as a single shared formula, mirroring the existing F543:F582
# shared-formula group that feeds the same column.
$ws.Range("F$startRow`:F$endRow").Formula = "=ABS(D$startRow-E$startRow)"

$ws.Range("R583:R589").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 565
